$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1646.125
$ws.Range("J17").Value = 1692.5333
$ws.Range("L17").Value = 5077.5999
$ws.Range("N17").Value = -5413.5999

$ws.Range("H18").Value = 2431.5
$ws.Range("I18").Value = 1797.8
$ws.Range("K18").Value = 1797.8
$ws.Range("M18").Value = -1513.8

$ws.Range("H44").Value = 21250
$ws.Range("J44").Value = 21250
$ws.Range("L44").Value = 21250
$ws.Range("N44").Value = -22174

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

$ws.Range("H135").Value = 2068.261
$ws.Range("I135").Value = 2061.111
$ws.Range("J135").Value = 2094
$ws.Range("K135").Value = 18549.999
$ws.Range("L135").Value = 18846
$ws.Range("M135").Value = -16014.999
$ws.Range("N135").Value = -23916

$ws.Range("H138").Value = 2475.6667
$ws.Range("J138").Value = 3078.6223
$ws.Range("L138").Value = 9235.866900000001
$ws.Range("N138").Value = -19515.8669

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29433262
$ws.Range("I32").Value = 45474936
$ws.Range("K32").Value = 45474936
$ws.Range("M32").Value = -45474649

$ws.Range("H61").Value = 41672384
$ws.Range("I61").Value = 35719600
$ws.Range("J61").Value = 62507132
$ws.Range("K61").Value = 35719600
$ws.Range("L61").Value = 62507132
$ws.Range("M61").Value = -35719388
$ws.Range("N61").Value = -62507556

$ws.Range("H74").Value = 10403535
$ws.Range("I74").Value = 14706844
$ws.Range("J74").Value = 1259004.9
$ws.Range("K74").Value = 14706844
$ws.Range("L74").Value = 1259004.9
$ws.Range("M74").Value = -14705970
$ws.Range("N74").Value = -1260752.9

$ws.Range("H77").Value = 10403535
$ws.Range("I77").Value = 14706844
$ws.Range("J77").Value = 1259004.9
$ws.Range("K77").Value = 73534220
$ws.Range("L77").Value = 6295024.5
$ws.Range("M77").Value = -73529852
$ws.Range("N77").Value = -6303760.5

$ws.Range("H136").Value = 41672384
$ws.Range("I136").Value = 35719600
$ws.Range("J136").Value = 62507132
$ws.Range("K136").Value = 107158800
$ws.Range("L136").Value = 187521396
$ws.Range("M136").Value = -107156250
$ws.Range("N136").Value = -187526496

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()

$ws.Range("H22").Value = 4191.2856
$ws.Range("I22").Value = 2889.8333
$ws.Range("K22").Value = 2889.8333
$ws.Range("M22").Value = -2716.8333

$ws.Range("H86").Value = 2091.4
$ws.Range("I86").Value = 2277.3157
$ws.Range("K86").Value = 2277.3157
$ws.Range("M86").Value = -1154.3157

$ws.Range("H89").Value = 2091.4
$ws.Range("I89").Value = 2277.3157
$ws.Range("K89").Value = 11386.5785
$ws.Range("M89").Value = -5770.5785

$ws.Range("H134").Value = 1001399.1
$ws.Range("I134").Value = 1495.375
$ws.Range("K134").Value = 4486.125
$ws.Range("M134").Value = -1951.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 524.1429000000001
$ws.Range("I22").Value = 533.9
$ws.Range("K22").Value = 533.9
$ws.Range("M22").Value = -183.9

$ws.Range("H31").Value = 408060.12
$ws.Range("I31").Value = 5986.357
$ws.Range("K31").Value = 5986.357
$ws.Range("M31").Value = -5691.357

$ws.Range("H34").Value = 408060.12
$ws.Range("I34").Value = 5986.357
$ws.Range("K34").Value = 5986.357
$ws.Range("M34").Value = -5784.357

$ws.Range("H57").Value = 20000
$ws.Range("J57").Value = 20000
$ws.Range("L57").Value = 20000
$ws.Range("N57").Value = -21120

$ws.Range("H68").Value = 61996.668
$ws.Range("J68").Value = 61996.668
$ws.Range("L68").Value = 61996.668
$ws.Range("N68").Value = -63494.668

$ws.Range("H71").Value = 61996.668
$ws.Range("J71").Value = 61996.668
$ws.Range("L71").Value = 185990.004
$ws.Range("N71").Value = -193478.004

$ws.Range("H94").Value = 5120.5
$ws.Range("J94").Value = 5358.75
$ws.Range("L94").Value = 5358.75
$ws.Range("N94").Value = -6260.75

$ws.Range("H107").Value = 1741.08
$ws.Range("I107").Value = 594.6429000000001
$ws.Range("K107").Value = 594.6429000000001
$ws.Range("M107").Value = 1325.3571

$ws.Range("H132").Value = 3583.7144
$ws.Range("I132").Value = 2848.0833
$ws.Range("J132").Value = 7997.5
$ws.Range("K132").Value = 8544.249899999999
$ws.Range("L132").Value = 23992.5
$ws.Range("M132").Value = -6014.249899999999
$ws.Range("N132").Value = -29052.5

$ws.Range("H134").Value = 3440.2144
$ws.Range("I134").Value = 2347.3635
$ws.Range("K134").Value = 7042.0905
$ws.Range("M134").Value = -4507.0905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 152.5
$ws.Range("I2").Value = 102.7
$ws.Range("K2").Value = 616.2
$ws.Range("M2").Value = -503.2

$ws.Range("H68").Value = 3166.5
$ws.Range("I68").Value = 999.5
$ws.Range("J68").Value = 4250
$ws.Range("K68").Value = 2998.5
$ws.Range("L68").Value = 12750
$ws.Range("M68").Value = -2187.5
$ws.Range("N68").Value = -14372

$ws.Range("H71").Value = 3166.5
$ws.Range("I71").Value = 999.5
$ws.Range("J71").Value = 4250
$ws.Range("K71").Value = 8995.5
$ws.Range("L71").Value = 38250
$ws.Range("M71").Value = -4939.5
$ws.Range("N71").Value = -46362

$ws.Range("H99").Value = 1008.6667
$ws.Range("I99").Value = 1008.6667
$ws.Range("K99").Value = 3026.0001
$ws.Range("M99").Value = -780.0001000000002

$ws.Range("H107").Value = 644.1786
$ws.Range("J107").Value = 994.6667
$ws.Range("L107").Value = 2984.0001
$ws.Range("N107").Value = -6824.0001

$ws.Range("H109").Value = 4000
$ws.Range("I109").Value = 4000
$ws.Range("K109").Value = 12000
$ws.Range("M109").Value = -10960

$ws.Range("H122").Value = 2886.389
$ws.Range("I122").Value = 371.42856
$ws.Range("J122").Value = 4486.8184
$ws.Range("K122").Value = 3342.85704
$ws.Range("L122").Value = 40381.3656
$ws.Range("M122").Value = -892.8570399999999
$ws.Range("N122").Value = -45281.3656

$ws.Range("H134").Value = 8634.643
$ws.Range("I134").Value = 1799.8
$ws.Range("J134").Value = 10120.479
$ws.Range("K134").Value = 5399.4
$ws.Range("L134").Value = 30361.437
$ws.Range("M134").Value = -329.3999999999996
$ws.Range("N134").Value = -40501.437

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 28000
$ws.Range("J48").Value = 28000
$ws.Range("L48").Value = 28000
$ws.Range("N48").Value = -28970

$ws.Range("H109").Value = 45141.5
$ws.Range("J109").Value = 45141.5
$ws.Range("L109").Value = 45141.5
$ws.Range("N109").Value = -47221.5

$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 71438590
$ws.Range("I132").Value = 142858670
$ws.Range("J132").Value = 18517.428
$ws.Range("K132").Value = 428576010
$ws.Range("L132").Value = 55552.284
$ws.Range("M132").Value = -428573480
$ws.Range("N132").Value = -60612.284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H25").Value = 14000

$ws.Range("H55").Value = 29412356
$ws.Range("I55").Value = 38462108
$ws.Range("K55").Value = 38462108
$ws.Range("M55").Value = -38461935

$ws.Range("H132").Value = 1264863.2
$ws.Range("I132").Value = 32634.666
$ws.Range("J132").Value = 2004200.4
$ws.Range("K132").Value = 97903.99800000001
$ws.Range("L132").Value = 6012601.199999999
$ws.Range("M132").Value = -95373.99800000001
$ws.Range("N132").Value = -6017661.199999999

$ws.Range("H136").Value = 99943.21000000001
$ws.Range("I136").Value = 15884.286
$ws.Range("K136").Value = 47652.858
$ws.Range("M136").Value = -45102.858

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 27601.4
$ws.Range("I15").Value = 17666.666
$ws.Range("J15").Value = 42503.5
$ws.Range("K15").Value = 17666.666
$ws.Range("L15").Value = 42503.5
$ws.Range("M15").Value = -17378.666
$ws.Range("N15").Value = -43079.5

$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H62").Value = 28578286
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 28578286
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H107").Value = 802.5333000000001
$ws.Range("I107").Value = 870.619
$ws.Range("K107").Value = 2611.857
$ws.Range("M107").Value = -691.857
